$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update cited_by_count for row 2 (Curcumin Supplementation paper): 9 -> 12
$ws.Range("M2").Value = 12

# 2) Rows 4-7 hold 4 publication records that need to be cyclically rotated:
#    old row4 -> new row5
#    old row5 -> new row6
#    old row6 -> new row7
#    old row7 -> new row4
# Use a scratch row (20) far outside the used range to stage the value of
# row 7 before it gets overwritten, then shift rows 6->7, 5->6, 4->5 and
# finally drop the stashed old-row-7 data into row 4.

# Stash old row 7 in scratch row 20
$ws.Range("A7:Q7").Copy($ws.Range("A20"))

# Shift old row 6 into row 7
$ws.Range("A6:Q6").Copy($ws.Range("A7"))

# Shift old row 5 into row 6
$ws.Range("A5:Q5").Copy($ws.Range("A6"))

# Shift old row 4 into row 5
$ws.Range("A4:Q4").Copy($ws.Range("A5"))

# Move stashed old row 7 into row 4
$ws.Range("A20:Q20").Copy($ws.Range("A4"))

# Clean up the scratch row so it doesn't linger in the saved workbook
$ws.Range("A20:Q20").ClearContents()

# 3) Fix author name in row 8: "Lisa B. Davidson" -> "Lisa Davidson"
$ws.Range("A8").Value = "Rocio Zapata Bustos, Dawn K. Coletta, Jean$([char]0x2010)Philippe Galons, Lisa Davidson, Paul Langlais, Janet L. Funk, Wayne T. Willis, Lawrence J. Mandarino"
